{"js": "// Replace the Finnish task-line text in the \"Laittomuuksia\" paragraph:\n//   \">>> Laittomuuksia (L\u00f6yd\u00e4 alakerta, Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 ja palauta t\u00e4m\u00e4 Koposelle)\"\n// becomes\n//   \">>> Laittomuuksia (Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 kemiavarastosta ja palauta t\u00e4m\u00e4 Koposelle)\"\n//\n// i.e. \"L\u00f6yd\u00e4 alakerta, \" is removed and \" kemiavarastosta\" is inserted right\n// after \"ep\u00e4ilytt\u00e4v\u00e4\u00e4\". The whole phrase is unique in the document, so we can\n// search for it directly and swap in the new wording (which also preserves\n// run-level formatting, since Word re-derives runs from the paragraph's\n// existing formatting when text is replaced in place).\n\nconst oldText =\n  \">>> Laittomuuksia (L\u00f6yd\u00e4 alakerta, Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 ja palauta t\u00e4m\u00e4 Koposelle)\";\nconst newText =\n  \">>> Laittomuuksia (Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 kemiavarastosta ja palauta t\u00e4m\u00e4 Koposelle)\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n} else {\n  // Fallback: perform the edit as a targeted removal + insertion in case the\n  // exact sentence couldn't be matched verbatim (e.g. whitespace drift).\n  const removeResults = body.search(\"L\u00f6yd\u00e4 alakerta, \", { matchCase: true });\n  removeResults.load(\"items\");\n  await context.sync();\n  removeResults.items.forEach((r) => r.insertText(\"\", Word.InsertLocation.replace));\n\n  const anchorResults = body.search(\"ep\u00e4ilytt\u00e4v\u00e4\u00e4\", { matchCase: true });\n  anchorResults.load(\"items\");\n  await context.sync();\n  if (anchorResults.items.length > 0) {\n    anchorResults.items[0].insertText(\" kemiavarastosta\", Word.InsertLocation.after);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the Finnish task-line text in the \"Laittomuuksia\" paragraph:\n#   \">>> Laittomuuksia (L\u00f6yd\u00e4 alakerta, Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 ja palauta t\u00e4m\u00e4 Koposelle)\"\n# becomes\n#   \">>> Laittomuuksia (Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 kemiavarastosta ja palauta t\u00e4m\u00e4 Koposelle)\"\n#\n# i.e. \"L\u00f6yd\u00e4 alakerta, \" is removed and \" kemiavarastosta\" is inserted right\n# after \"ep\u00e4ilytt\u00e4v\u00e4\u00e4\". The phrase is unique in the document, so Find/Replace\n# on the whole sentence is the most robust way to make the swap.\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$oldText = \">>> Laittomuuksia (L\u00f6yd\u00e4 alakerta, Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 ja palauta t\u00e4m\u00e4 Koposelle)\"\n$newText = \">>> Laittomuuksia (Etsi jotain ep\u00e4ilytt\u00e4v\u00e4\u00e4 kemiavarastosta ja palauta t\u00e4m\u00e4 Koposelle)\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne)\n\nif (-not $found) {\n    # Fallback: perform the edit as a targeted removal + insertion in case the\n    # exact sentence couldn't be matched verbatim (e.g. whitespace drift).\n    $findRemove = $d.Content.Find\n    $findRemove.Text = \"L\u00f6yd\u00e4 alakerta, \"\n    $findRemove.Replacement.Text = \"\"\n    $findRemove.Execute($findRemove.Text, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $findRemove.Replacement.Text, $wdReplaceOne)\n\n    $findAnchor = $d.Content.Find\n    $findAnchor.Text = \"ep\u00e4ilytt\u00e4v\u00e4\u00e4\"\n    $anchorFound = $findAnchor.Execute($findAnchor.Text, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, \"\", 0)\n    if ($anchorFound) {\n        $hit = $findAnchor.Parent\n        $hit.Collapse(0)\n        $hit.InsertAfter(\" kemiavarastosta\")\n    }\n}\n"}
